$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row - add "Total Amount" header in C1
$ws.Range("C1").Value = "Total Amount"

# Update the date portion of the file paths in column A (06-30-2022 -> 07-01-2022)
$ws.Range("A2").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\07-01-2022\Axis\CN0138305133_230522_18_10_29.pdf"
$ws.Range("A3").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\07-01-2022\Axis\CN0138305152_230522_18_10_29.pdf"
$ws.Range("A4").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\07-01-2022\Axis\CN0138305153_230522_18_10_29.pdf"
$ws.Range("A5").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\07-01-2022\Axis\CR0006093473_230522_18_10_29.pdf"
$ws.Range("A6").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\07-01-2022\Axis\CR0006093474_230522_18_10_29.pdf"

# Add Total Amount values in column C
$ws.Range("C2").Value = "INR 1,64,250.00"
$ws.Range("C3").Value = "INR 35,231.78"
$ws.Range("C4").Value = "INR 32,928.00"
$ws.Range("C5").Value = "INR 17,99,020.40"
$ws.Range("C6").Value = "INR 15,16,082.00"

# Set column C width to match bestFit width used elsewhere (approx 15)
$ws.Columns.Item(3).ColumnWidth = 14.1667

# Select column A fully (A1:A1048576) to match the sheetView selection change
$ws.Range("A1:A1048576").Select()
